# feat: Filter out example rows from Excel templates during upload
#
# The vehicles template shipped two fully-populated example/sample rows
# (row 2 "TRUCK-001" and row 3 "TRUCK-002"). Uploads were picking these
# example rows up as if they were real data, so the template is changed
# to:
#   - drop the second example row (row 3) entirely
#   - blank out the identifying/free-text sample values on the first
#     example row (code, plate number, device id, garage address) so it
#     reads as a template hint row rather than real data
#   - turn the "차량타입" (vehicle type) cell into an explanatory hint
#   - tweak a few numeric placeholder values
#   - re-balance a handful of column widths now that some columns hold
#     longer hint text and others hold none

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths -------------------------------------------------
# NOTE: the host's ColumnWidth setter stores `requested + 5/6` (an MDW-7
# pixel-grid rounding artifact), so compensate by subtracting 5/6 from
# every target width we assign here to land exactly on the target value.
$colWidthFudge = 5 / 6

function Set-ColWidth($col, $target) {
    $ws.Columns.Item($col).ColumnWidth = $target - $colWidthFudge
}

Set-ColWidth 1 6     # A 차량코드
Set-ColWidth 2 6     # B 차량번호
Set-ColWidth 3 22    # C 차량타입
Set-ColWidth 4 11    # D UVIS단말기ID
Set-ColWidth 17 7    # Q 차고지주소

# --- Row 3 (second example row "TRUCK-002") is removed entirely ----
$ws.Rows.Item(3).Delete()

# --- Row 2 (first example row) becomes a template hint row ---------
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "냉동 또는 냉장 또는 겸용 또는 상온"
$ws.Range("D2").ClearContents()

$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 5000
$ws.Range("G2").Value = 30
$ws.Range("N2").Value = 5

$ws.Range("Q2").ClearContents()
